# Edit: 2020-07-13 commit
# 1) Slide 6 table: switch the table's style from the custom "Table_0"
#    style to the built-in PowerPoint table style
#    {674D182B-7514-4622-945F-3056F27C638B}.
# 2) Presentation design: the deck's theme (used by the slide master and
#    all slides/layouts) changes its 12 theme colors from the "Integral"
#    palette to the standard "Office Theme" palette.

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 6 -------------------------------------------
$slide = $p.Slides.Item(6)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{674D182B-7514-4622-945F-3056F27C638B}", $true)
    }
}

# --- 2) Theme colors: Integral -> Office Theme ---------------------------
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

$tcs.Item(1).RGB  = 0          # dk1      000000
$tcs.Item(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      44546A
$tcs.Item(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  4472C4
$tcs.Item(10).RGB = 4697456    # accent6  70AD47
$tcs.Item(11).RGB = 12673797   # hlink    0563C1
$tcs.Item(12).RGB = 7491477    # folHlink 954F72
